$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 869.2308
$ws.Range("I11").Value = 869.2308
$ws.Range("K11").Value = 869.2308
$ws.Range("M11").Value = -729.2308
$ws.Range("H53").Value = 4278.1875
$ws.Range("I53").Value = 3712.2
$ws.Range("J53").Value = 5221.5
$ws.Range("K53").Value = 3712.2
$ws.Range("L53").Value = 5221.5
$ws.Range("M53").Value = -3075.2
$ws.Range("N53").Value = -6495.5
$ws.Range("H86").Value = 99912190
$ws.Range("J86").Value = 18521850
$ws.Range("L86").Value = 18521850
$ws.Range("N86").Value = -18524096
$ws.Range("H89").Value = 99912190
$ws.Range("J89").Value = 18521850
$ws.Range("L89").Value = 92609250
$ws.Range("N89").Value = -92620482
$ws.Range("H112").Value = 3955.9375
$ws.Range("J112").Value = 3955.9375
$ws.Range("L112").Value = 11867.8125
$ws.Range("N112").Value = -14083.8125
$ws.Range("H132").Value = 770.63635
$ws.Range("I132").Value = 770.63635
$ws.Range("K132").Value = 2311.90905
$ws.Range("M132").Value = 218.0909499999998
$ws.Range("H137").Value = 4581.7
$ws.Range("I137").Value = 5799.4
$ws.Range("K137").Value = 17398.2
$ws.Range("M137").Value = -14848.2
$ws.Range("H138").Value = 7641.2256
$ws.Range("I138").Value = 3984.8667
$ws.Range("J138").Value = 11069.0625
$ws.Range("K138").Value = 11954.6001
$ws.Range("L138").Value = 33207.1875
$ws.Range("M138").Value = -6814.6001
$ws.Range("N138").Value = -43487.1875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2010311
$ws.Range("I32").Value = 2030056.1
$ws.Range("K32").Value = 2030056.1
$ws.Range("M32").Value = -2029769.1
$ws.Range("H61").Value = 35721532
$ws.Range("I61").Value = 3774.6667
$ws.Range("K61").Value = 3774.6667
$ws.Range("M61").Value = -3562.6667
$ws.Range("H122").Value = 2392
$ws.Range("I122").Value = 2392
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 7176
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -4726
$ws.Range("H125").Value = 55569.92
$ws.Range("J125").Value = 55569.92
$ws.Range("L125").Value = 55569.92
$ws.Range("N125").Value = -65409.92
$ws.Range("H132").Value = 4955.6377
$ws.Range("I132").Value = 2690.6875
$ws.Range("J132").Value = 7743.269
$ws.Range("K132").Value = 8072.0625
$ws.Range("L132").Value = 23229.807
$ws.Range("M132").Value = -5542.0625
$ws.Range("N132").Value = -28289.807
$ws.Range("H136").Value = 35721532
$ws.Range("I136").Value = 3774.6667
$ws.Range("K136").Value = 11324.0001
$ws.Range("M136").Value = -8774.000100000001
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 52021.855
$ws.Range("I86").Value = 70859.2
$ws.Range("J86").Value = 4928.5
$ws.Range("K86").Value = 70859.2
$ws.Range("L86").Value = 4928.5
$ws.Range("M86").Value = -69736.2
$ws.Range("N86").Value = -7174.5
$ws.Range("H89").Value = 52021.855
$ws.Range("I89").Value = 70859.2
$ws.Range("J89").Value = 4928.5
$ws.Range("K89").Value = 354296
$ws.Range("L89").Value = 24642.5
$ws.Range("M89").Value = -348680
$ws.Range("N89").Value = -35874.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 25790.51
$ws.Range("I31").Value = 2926.3704
$ws.Range("J31").Value = 51512.668
$ws.Range("K31").Value = 2926.3704
$ws.Range("L31").Value = 51512.668
$ws.Range("M31").Value = -2631.3704
$ws.Range("N31").Value = -52102.668
$ws.Range("H34").Value = 25790.51
$ws.Range("I34").Value = 2926.3704
$ws.Range("J34").Value = 51512.668
$ws.Range("K34").Value = 2926.3704
$ws.Range("L34").Value = 51512.668
$ws.Range("M34").Value = -2724.3704
$ws.Range("N34").Value = -51916.668
$ws.Range("H99").Value = 7132.6
$ws.Range("I99").Value = 6144
$ws.Range("K99").Value = 6144
$ws.Range("M99").Value = -4646
$ws.Range("H122").Value = 4103.0557
$ws.Range("J122").Value = 5656
$ws.Range("L122").Value = 16968
$ws.Range("N122").Value = -21868
$ws.Range("H124").Value = 71999.664
$ws.Range("J124").Value = 71999.664
$ws.Range("L124").Value = 71999.664
$ws.Range("N124").Value = -76909.664
$ws.Range("H126").Value = 7132.6
$ws.Range("I126").Value = 6144
$ws.Range("K126").Value = 18432
$ws.Range("M126").Value = -15962
$ws.Range("H132").Value = 7830.968
$ws.Range("I132").Value = 5881.6
$ws.Range("J132").Value = 9658.5
$ws.Range("K132").Value = 17644.8
$ws.Range("L132").Value = 28975.5
$ws.Range("M132").Value = -15114.8
$ws.Range("N132").Value = -34035.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 7900
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 7900
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 23700
$ws.Range("N94").Value = -25052
$ws.Range("H113").Value = 3655.875
$ws.Range("J113").Value = 3655.875
$ws.Range("L113").Value = 10967.625
$ws.Range("N113").Value = -15307.625
$ws.Range("H122").Value = 2016802.8
$ws.Range("J122").Value = 1001602.8
$ws.Range("L122").Value = 9014425.200000001
$ws.Range("N122").Value = -9019325.200000001
$ws.Range("H131").Value = 40588.883
$ws.Range("J131").Value = 68990.92999999999
$ws.Range("L131").Value = 206972.79
$ws.Range("N131").Value = -217052.79
$ws.Range("H132").Value = 6980.4062
$ws.Range("I132").Value = 5977
$ws.Range("K132").Value = 53793
$ws.Range("M132").Value = -51263
$ws.Range("M94").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 116012.664
$ws.Range("J80").Value = 172519.17
$ws.Range("L80").Value = 172519.17
$ws.Range("N80").Value = -174515.17
$ws.Range("H83").Value = 116012.664
$ws.Range("J83").Value = 172519.17
$ws.Range("L83").Value = 862595.8500000001
$ws.Range("N83").Value = -872579.8500000001
$ws.Range("H132").Value = 5162.475
$ws.Range("I132").Value = 3381.5454
$ws.Range("K132").Value = 10144.6362
$ws.Range("M132").Value = -7614.636200000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6300.1
$ws.Range("I7").Value = 6166.8335
$ws.Range("K7").Value = 6166.8335
$ws.Range("M7").Value = -6054.8335
$ws.Range("H40").Value = 8699.362999999999
$ws.Range("I40").Value = 5596.5
$ws.Range("J40").Value = 9388.888999999999
$ws.Range("K40").Value = 5596.5
$ws.Range("L40").Value = 9388.888999999999
$ws.Range("M40").Value = -5460.5
$ws.Range("N40").Value = -9660.888999999999
$ws.Range("H46").Value = 1872.3
$ws.Range("I46").Value = 975.4286
$ws.Range("K46").Value = 975.4286
$ws.Range("M46").Value = -787.4286
$ws.Range("H122").Value = 5515.548
$ws.Range("I122").Value = 4909.0835
$ws.Range("J122").Value = 5758.1333
$ws.Range("K122").Value = 14727.2505
$ws.Range("L122").Value = 17274.3999
$ws.Range("M122").Value = -12277.2505
$ws.Range("N122").Value = -22174.3999
$ws.Range("H126").Value = 6300.1
$ws.Range("I126").Value = 6166.8335
$ws.Range("K126").Value = 18500.5005
$ws.Range("M126").Value = -16030.5005
$ws.Range("H132").Value = 11912747
$ws.Range("I132").Value = 17862298
$ws.Range("K132").Value = 53586894
$ws.Range("M132").Value = -53584364

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 500500000
$ws.Range("I14").Value = 500500000
$ws.Range("K14").Value = 500500000
$ws.Range("M14").Value = -500499832
$ws.Range("H75").Value = 45372.668
$ws.Range("I75").Value = 45372.668
$ws.Range("K75").Value = 45372.668
$ws.Range("M75").Value = -44436.668
$ws.Range("H78").Value = 45372.668
$ws.Range("I78").Value = 45372.668
$ws.Range("K78").Value = 136118.004
$ws.Range("M78").Value = -131438.004
$ws.Range("H100").Value = 564.94446
$ws.Range("I100").Value = 446.22223
$ws.Range("K100").Value = 892.44446
$ws.Range("M100").Value = -351.44446
$ws.Range("H122").Value = 109971.69
$ws.Range("J122").Value = 8123.077
$ws.Range("L122").Value = 24369.231
$ws.Range("N122").Value = -29269.231
$ws.Range("H125").Value = 49806.668
$ws.Range("J125").Value = 49806.668
$ws.Range("L125").Value = 49806.668
$ws.Range("N125").Value = -59646.668
$ws.Range("H132").Value = 8190.591
$ws.Range("I132").Value = 10273.52
$ws.Range("K132").Value = 30820.56
$ws.Range("M132").Value = -28290.56
$ws.Range("H135").Value = 63703.715
$ws.Range("J135").Value = 63703.715
$ws.Range("L135").Value = 63703.715
$ws.Range("N135").Value = -73843.715
$ws.Range("H140").Value = 149000
$ws.Range("J140").Value = 149000
$ws.Range("L140").Value = 149000
$ws.Range("N140").Value = -159360
